$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (pushes the existing rows 8-26 down to 9-27)
$ws.Rows.Item(8).Insert()

# Copy the formatting of the row below (the row that used to be row 8, now row 9)
# into the freshly inserted row 8 so its cell styles match the rest of the table.
$ws.Range("A9:E9").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new task row: "Mise en fonctionnement de l'API"
$ws.Range("A8").Value = "Mise en fonctionnement de l'API"
$ws.Range("B8").Value = "Pierre"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "Android Studio, VS Code"
$ws.Range("E8").Value = "Pas abouti"

# The existing "Réalisation de l'API" row (now row 7) also gained VS Code as a tool
$ws.Range("D7").Value = "Android Studio, VS Code"

# Recalculate so the SUMIF/SUM totals in column G reflect the new row
$excel.Calculate()

# Update the view's current selection to match where the author ended up
[void]$ws.Range("E9").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
